$wb = $excel.ActiveWorkbook

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 12:05:41"

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 12:05:31"
$wsZhCn.Range("K4").Value = "2016-08-31 12:06:25"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-31 12:05:41"
$wsDeDe.Range("K4").Value = "2016-08-31 12:06:42"
